# Add a new row "HubSante.etape / Etape_Message / ENUM-Etape_Message / ENUM"
# to the nomenclature table, right after the "ISO 3166-ISO3166-2" row and
# before the "HubSante.statutVecteur" row.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row that currently immediately follows the "ISO 3166" row
# (i.e. the row whose first cell reads "HubSante.statutVecteur"), so the
# new row can be inserted right before it.
$targetRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 1).Range.Text
    if ($cellText -like "HubSante.statutVecteur*") {
        $targetRow = $t.Rows.Item($i)
        break
    }
}

$newRow = $t.Rows.Add($targetRow)
$newRow.Cells.Item(1).Range.Text = "HubSante.etape"
$newRow.Cells.Item(2).Range.Text = "Etape_Message"
$newRow.Cells.Item(3).Range.Text = "ENUM-Etape_Message"
$newRow.Cells.Item(4).Range.Text = "ENUM"
